$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1), columns T:X ---
$ws.Range("S1").Copy() | Out-Null
$ws.Range("T1:X1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("T1").Value = " Tarjetas credito vigentes otros"
$ws.Range("U1").Value = "Numero de operaciones realizadas con tarjetas de crédito"
$ws.Range("V1").Value = "Tarjetas vigentes"
$ws.Range("W1").Value = "Tarjetas vigentes VISA"
$ws.Range("X1").Value = "Tarjetas vigentes MASTERCARD"

# --- Row 2 (SMAPE) new values ---
$ws.Range("T2").Value = 0.1515117755174033
$ws.Range("U2").Value = 0.04893548824463104
$ws.Range("V2").Value = 0.03035588416534699
$ws.Range("W2").Value = 0.04510072663549621
$ws.Range("X2").Value = 0.02048521563172385

# --- Row 3 (MAE) new values ---
$ws.Range("T3").Value = 67254.58748056635
$ws.Range("U3").Value = 402297.5488382598
$ws.Range("V3").Value = 12294202.13725019
$ws.Range("W3").Value = 2096838.884939871
$ws.Range("X3").Value = 594844.9597643962

# --- Row 4 (MASE) new values ---
$ws.Range("T4").Value = 0.3678116291786256
$ws.Range("U4").Value = 0.08291884173554066
$ws.Range("V4").Value = 0.05813484878731042
$ws.Range("W4").Value = 0.08886124252627792
$ws.Range("X4").Value = 0.04251739920413474

# --- New row 5: Confiabilidad 80% ---
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A5").Value = "Confiabilidad 80%"

$ws.Range("B5").Value = 0.1666666666666667
$ws.Range("C5").Value = 0.1666666666666667
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0.3333333333333333
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.1666666666666667
$ws.Range("T5").Value = 0.1666666666666667
$ws.Range("U5").Value = 0.8333333333333334
$ws.Range("V5").Value = 1
$ws.Range("W5").Value = 1
$ws.Range("X5").Value = 1
